$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ML Project Checklist")

# Mark "Done" column (C) for rows 11-16 with "Y" and rows 17-20 with "y"
$ws.Range("C11").Value = "Y"
$ws.Range("C12").Value = "Y"
$ws.Range("C13").Value = "Y"
$ws.Range("C14").Value = "Y"
$ws.Range("C15").Value = "Y"
$ws.Range("C16").Value = "Y"
$ws.Range("C17").Value = "y"
$ws.Range("C18").Value = "y"
$ws.Range("C19").Value = "y"
$ws.Range("C20").Value = "y"

# Update the view: scroll so row 15 is the top row, and select C21
$ws.Range("C21").Select()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
